$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 200006500
$ws.Range("J64").Value = 250005870
$ws.Range("L64").Value = 250005870
$ws.Range("N64").Value = -250006366
$ws.Range("H67").Value = 200006500
$ws.Range("J67").Value = 250005870
$ws.Range("L67").Value = 250005870
$ws.Range("N67").Value = -250007586
$ws.Range("H80").Value = 419.41666
$ws.Range("I80").Value = 535.5714
$ws.Range("J80").Value = 256.8
$ws.Range("K80").Value = 1606.7142
$ws.Range("L80").Value = 770.4000000000001
$ws.Range("M80").Value = -608.7142000000001
$ws.Range("N80").Value = -2766.4
$ws.Range("H83").Value = 419.41666
$ws.Range("I83").Value = 535.5714
$ws.Range("J83").Value = 256.8
$ws.Range("K83").Value = 4820.1426
$ws.Range("L83").Value = 2311.2
$ws.Range("M83").Value = 171.8573999999999
$ws.Range("N83").Value = -12295.2
$ws.Range("H98").Value = 4219.636
$ws.Range("I98").Value = 4219.636
$ws.Range("K98").Value = 4219.636
$ws.Range("M98").Value = -2721.636
$ws.Range("H100").Value = 2259.111
$ws.Range("I100").Value = 1291.625
$ws.Range("K100").Value = 1291.625
$ws.Range("M100").Value = -750.625
$ws.Range("H101").Value = 375
$ws.Range("I101").Value = 375
$ws.Range("K101").Value = 1125
$ws.Range("M101").Value = 497
$ws.Range("H113").Value = 3708.9092
$ws.Range("I113").Value = 2800
$ws.Range("J113").Value = 3799.8
$ws.Range("K113").Value = 2800
$ws.Range("L113").Value = 3799.8
$ws.Range("M113").Value = 454
$ws.Range("N113").Value = -10307.8
$ws.Range("H122").Value = 4219.636
$ws.Range("I122").Value = 4219.636
$ws.Range("K122").Value = 12658.908
$ws.Range("M122").Value = -10208.908
$ws.Range("H125").Value = 1833.2222
$ws.Range("I125").Value = 999
$ws.Range("J125").Value = 1937.5
$ws.Range("K125").Value = 8991
$ws.Range("L125").Value = 17437.5
$ws.Range("M125").Value = -6531
$ws.Range("N125").Value = -22357.5
$ws.Range("H135").Value = 2552.7144
$ws.Range("I135").Value = 1560.4
$ws.Range("K135").Value = 14043.6
$ws.Range("M135").Value = -11508.6
$ws.Range("H137").Value = 2506166.5
$ws.Range("I137").Value = 2941960.5
$ws.Range("K137").Value = 8825881.5
$ws.Range("M137").Value = -8823331.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7429.7
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 7429.7
$ws.Range("K32").Value = 0
$ws.Range("L32").ClearContents()
$ws.Range("M32").Value = 7429.7
$ws.Range("N32").Value = -8003.7
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 89333.336
$ws.Range("J132").Value = 89333.336
$ws.Range("L132").Value = 89333.336
$ws.Range("N132").Value = -99453.336
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1502.8
$ws.Range("I16").Value = 1339.5416
$ws.Range("K16").Value = 1339.5416
$ws.Range("M16").Value = -1052.5416
$ws.Range("H58").Value = 2883.5652
$ws.Range("I58").Value = 2108.5715
$ws.Range("J58").Value = 4089.111
$ws.Range("K58").Value = 2108.5715
$ws.Range("L58").Value = 4089.111
$ws.Range("M58").Value = -1905.5715
$ws.Range("N58").Value = -4495.111
$ws.Range("H99").Value = 6559.6
$ws.Range("J99").Value = 7599.6665
$ws.Range("L99").Value = 7599.6665
$ws.Range("N99").Value = -10595.6665
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").ClearContents()
$ws.Range("N111").Value = 0
$ws.Range("H113").Value = 1502.8
$ws.Range("I113").Value = 1339.5416
$ws.Range("K113").Value = 1339.5416
$ws.Range("M113").Value = 830.4584
$ws.Range("H126").Value = 6559.6
$ws.Range("J126").Value = 7599.6665
$ws.Range("L126").Value = 22798.9995
$ws.Range("N126").Value = -27738.9995
$ws.Range("H132").Value = 2710.1
$ws.Range("J132").Value = 4000
$ws.Range("L132").Value = 12000
$ws.Range("N132").Value = -17060
$ws.Range("H134").Value = 2022.1875
$ws.Range("I134").Value = 1958.3871
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 5875.1613
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -3340.1613
$ws.Range("N134").Value = -17070
$ws.Range("H136").Value = 2883.5652
$ws.Range("I136").Value = 2108.5715
$ws.Range("J136").Value = 4089.111
$ws.Range("K136").Value = 6325.7145
$ws.Range("L136").Value = 12267.333
$ws.Range("M136").Value = -3775.7145
$ws.Range("N136").Value = -17367.333
$ws.Range("H140").Value = 94999.25
$ws.Range("J140").Value = 98713.42999999999
$ws.Range("L140").Value = 98713.42999999999
$ws.Range("N140").Value = -109073.43
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 751
$ws.Range("I134").Value = 751
$ws.Range("K134").Value = 2253
$ws.Range("M134").Value = 2817
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 5333.25
$ws.Range("I113").Value = 4533.9
$ws.Range("K113").Value = 4533.9
$ws.Range("M113").Value = -2363.9
$ws.Range("H132").Value = 3095.7144
$ws.Range("I132").Value = 2472.7778
$ws.Range("K132").Value = 7418.3334
$ws.Range("M132").Value = -4888.3334
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2162.818
$ws.Range("I7").Value = 2257.4
$ws.Range("J7").Value = 2084
$ws.Range("K7").Value = 2257.4
$ws.Range("L7").Value = 2084
$ws.Range("M7").Value = -2145.4
$ws.Range("N7").Value = -2308
$ws.Range("H22").Value = 800.25
$ws.Range("I22").Value = 725.5
$ws.Range("J22").Value = 875
$ws.Range("K22").Value = 725.5
$ws.Range("L22").Value = 875
$ws.Range("M22").Value = -430.5
$ws.Range("N22").Value = -1465
$ws.Range("H27").Value = 800.25
$ws.Range("I27").Value = 725.5
$ws.Range("J27").Value = 875
$ws.Range("K27").Value = 725.5
$ws.Range("L27").Value = 875
$ws.Range("M27").Value = -618.5
$ws.Range("N27").Value = -1089
$ws.Range("H46").Value = 1399.3334
$ws.Range("I46").Value = 1399
$ws.Range("J46").Value = 1400
$ws.Range("K46").Value = 1399
$ws.Range("L46").Value = 1400
$ws.Range("M46").Value = -1211
$ws.Range("N46").Value = -1776
$ws.Range("H126").Value = 2162.818
$ws.Range("I126").Value = 2257.4
$ws.Range("J126").Value = 2084
$ws.Range("K126").Value = 6772.200000000001
$ws.Range("L126").Value = 6252
$ws.Range("M126").Value = -4302.200000000001
$ws.Range("N126").Value = -11192
$ws.Range("H132").Value = 6773.143
$ws.Range("I132").Value = 3793.5557
$ws.Range("K132").Value = 11380.6671
$ws.Range("M132").Value = -8850.667099999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1027.2
$ws.Range("I126").Value = 1027.2
$ws.Range("K126").Value = 3081.6
$ws.Range("M126").Value = -611.6000000000004
$ws.Range("H136").Value = 83338750
$ws.Range("I136").Value = 125001250
$ws.Range("J136").Value = 13750
$ws.Range("K136").Value = 375003750
$ws.Range("L136").Value = 41250
$ws.Range("M136").Value = -375001200
$ws.Range("N136").Value = -46350
